# Added PN100 - Rear PTFE Tube Guide to the printed-parts BOM table.
# The table is sorted by Number (column A), and PN100 sorts right after
# PN099 (row 74) / before PN102 (old row 75), so the new row is inserted
# at sheet row 75, pushing the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("000 - Printed Parts")
$lo = $ws.ListObjects.Item("Table1")

# Insert a new, blank row at row 75 - this shifts the existing PN102.. rows
# (and everything below them) down by one row, preserving their data.
$ws.Rows("75:75").Insert()

# Grow the table / autofilter range by one row to cover the new row.
$lo.Resize($ws.Range("A1:I96"))

# Fill in the new part row: PN100, XY / Wiring, no chirality (N), name,
# material ABS, qty 1, STL filename. (No Note for this row.)
$ws.Range("A75").Value2 = "PN100"
$ws.Range("B75").Value2 = "XY"
$ws.Range("C75").Value2 = "Wiring"
$ws.Range("D75").Value2 = "N"
$ws.Range("E75").Value2 = "Rear PTFE Tube Guide"
$ws.Range("F75").Value2 = "ABS"
$ws.Range("G75").Value2 = 1
$ws.Range("I75").Value2 = "100 - XY - Wiring - Rear PTFE Tube Guide.stl"

# Restore the saved view/selection state (scrolled down near the new row,
# with I75 selected).
$ws.Range("A70").Select()
$ws.Range("I75").Select()
